$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1100.1538
$ws.Range("J17").Value = 981.1905
$ws.Range("L17").Value = 2943.5715
$ws.Range("N17").Value = -3279.5715

# Row 28
$ws.Range("H28").Value = 520.2
$ws.Range("I28").Value = 401.75
$ws.Range("K28").Value = 401.75
$ws.Range("M28").Value = 83.25

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 51
$ws.Range("H51").Value = 100030000
$ws.Range("I51").Value = 250027500
$ws.Range("J51").Value = 31666.666
$ws.Range("K51").Value = 250027500
$ws.Range("L51").Value = 31666.666
$ws.Range("M51").Value = -250027016
$ws.Range("N51").Value = -32634.666

# Row 58
$ws.Range("H58").Value = 263.5
$ws.Range("I58").Value = 263.5
$ws.Range("K58").Value = 790.5
$ws.Range("M58").Value = -640.5

# Row 111
$ws.Range("H111").Value = 3790.4614
$ws.Range("J111").Value = 1712.4286
$ws.Range("L111").Value = 5137.2858
$ws.Range("N111").Value = -11271.2858

# Row 132
$ws.Range("H132").Value = 5845.2173
$ws.Range("I132").Value = 5845.2173
$ws.Range("K132").Value = 17535.6519
$ws.Range("M132").Value = -15005.6519

# Row 135
$ws.Range("H135").Value = 815.125
$ws.Range("I135").Value = 385.30768
$ws.Range("K135").Value = 3467.76912
$ws.Range("M135").Value = -932.7691199999999

# Row 137
$ws.Range("H137").Value = 13644.272
$ws.Range("J137").Value = 26800.25
$ws.Range("L137").Value = 80400.75
$ws.Range("N137").Value = -85500.75

# Row 138
$ws.Range("H138").Value = 1323.64
$ws.Range("I138").Value = 718.619
$ws.Range("K138").Value = 2155.857
$ws.Range("M138").Value = 2984.143

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 90
$ws.Range("K5").Value = 90
$ws.Range("M5").Value = 22

# Row 61
$ws.Range("H61").Value = 2630.76
$ws.Range("I61").Value = 1523.6666
$ws.Range("K61").Value = 1523.6666
$ws.Range("M61").Value = -1311.6666

# Row 74
$ws.Range("H74").Value = 233347.2
$ws.Range("I74").Value = 293953.47
$ws.Range("K74").Value = 293953.47
$ws.Range("M74").Value = -293079.47

# Row 77
$ws.Range("H77").Value = 233347.2
$ws.Range("I77").Value = 293953.47
$ws.Range("K77").Value = 1469767.35
$ws.Range("M77").Value = -1465399.35

# Row 122
$ws.Range("H122").Value = 3799.8572
$ws.Range("I122").Value = 4066.5833
$ws.Range("K122").Value = 12199.7499
$ws.Range("M122").Value = -9749.749899999999

# Row 132
$ws.Range("H132").Value = 2453.4138
$ws.Range("I132").Value = 2297.9167
$ws.Range("K132").Value = 6893.750100000001
$ws.Range("M132").Value = -4363.750100000001

# Row 136
$ws.Range("H136").Value = 2630.76
$ws.Range("I136").Value = 1523.6666
$ws.Range("K136").Value = 4570.9998
$ws.Range("M136").Value = -2020.9998

# Row 139
$ws.Range("H139").Value = 69749.164
$ws.Range("I139").Value = 68000
$ws.Range("J139").Value = 69908.17999999999
$ws.Range("K139").Value = 68000
$ws.Range("L139").Value = 69908.17999999999
$ws.Range("M139").Value = -62860
$ws.Range("N139").Value = -80188.17999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 90
$ws.Range("M4").Value = 25

# Row 64
$ws.Range("H64").Value = 1522.4615
$ws.Range("J64").Value = 1798.75
$ws.Range("L64").Value = 1798.75
$ws.Range("N64").Value = -2248.75

# Row 67
$ws.Range("H67").Value = 1522.4615
$ws.Range("J67").Value = 1798.75
$ws.Range("L67").Value = 1798.75
$ws.Range("N67").Value = -3358.75

# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 94
$ws.Range("H94").Value = 83338060
$ws.Range("I94").Value = 95243360
$ws.Range("K94").Value = 95243360
$ws.Range("M94").Value = -95242909

# Row 134
$ws.Range("H134").Value = 2152.2942
$ws.Range("I134").Value = 1737.6154
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 5212.8462
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -2677.8462
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3680028
$ws.Range("J31").Value = 17863436
$ws.Range("L31").Value = 17863436
$ws.Range("N31").Value = -17864026

# Row 34
$ws.Range("H34").Value = 3680028
$ws.Range("J34").Value = 17863436
$ws.Range("L34").Value = 17863436
$ws.Range("N34").Value = -17863840

# Row 58
$ws.Range("H58").Value = 2418.1
$ws.Range("I58").Value = 1636.2
$ws.Range("K58").Value = 1636.2
$ws.Range("M58").Value = -1433.2

# Row 94
$ws.Range("H94").Value = 693.9091
$ws.Range("J94").Value = 749.4286
$ws.Range("L94").Value = 749.4286
$ws.Range("N94").Value = -1651.4286

# Row 99
$ws.Range("H99").Value = 8000
$ws.Range("J99").Value = 7000
$ws.Range("L99").Value = 7000
$ws.Range("N99").Value = -9996

# Row 107
$ws.Range("H107").Value = 2273781.8
$ws.Range("I107").Value = 3334130.2
$ws.Range("J107").Value = 1606.5714
$ws.Range("K107").Value = 3334130.2
$ws.Range("L107").Value = 1606.5714
$ws.Range("M107").Value = -3332210.2
$ws.Range("N107").Value = -5446.5714

# Row 122
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19900

# Row 126
$ws.Range("H126").Value = 8000
$ws.Range("J126").Value = 7000
$ws.Range("L126").Value = 21000
$ws.Range("N126").Value = -25940

# Row 132
$ws.Range("H132").Value = 3418.439
$ws.Range("I132").Value = 2948.5518
$ws.Range("J132").Value = 4554
$ws.Range("K132").Value = 8845.6554
$ws.Range("L132").Value = 13662
$ws.Range("M132").Value = -6315.6554
$ws.Range("N132").Value = -18722

# Row 134
$ws.Range("H134").Value = 4486.56
$ws.Range("I134").Value = 4538.636
$ws.Range("K134").Value = 13615.908
$ws.Range("M134").Value = -11080.908

# Row 136
$ws.Range("H136").Value = 2418.1
$ws.Range("I136").Value = 1636.2
$ws.Range("K136").Value = 4908.6
$ws.Range("M136").Value = -2358.6

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 1998.25
$ws.Range("J92").Value = 1998
$ws.Range("L92").Value = 5994
$ws.Range("N92").Value = -8490

# Row 103
$ws.Range("H103").Value = 2731.2
$ws.Range("J103").Value = 2731.2
$ws.Range("L103").Value = 8193.599999999999
$ws.Range("N103").Value = -9951.599999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1701.8379
$ws.Range("I102").Value = 1240.1538
$ws.Range("K102").Value = 1240.1538
$ws.Range("M102").Value = 381.8462

# Row 126
$ws.Range("H126").Value = 8968.15
$ws.Range("I126").Value = 2146.8333
$ws.Range("K126").Value = 6440.499899999999
$ws.Range("M126").Value = -3970.499899999999

# Row 132
$ws.Range("H132").Value = 2310
$ws.Range("I132").Value = 1908.2307
$ws.Range("J132").Value = 2658.2
$ws.Range("K132").Value = 5724.6921
$ws.Range("L132").Value = 7974.599999999999
$ws.Range("M132").Value = -3194.6921
$ws.Range("N132").Value = -13034.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1180.9375
$ws.Range("I22").Value = 1247.4166
$ws.Range("J22").Value = 981.5
$ws.Range("K22").Value = 1247.4166
$ws.Range("L22").Value = 981.5
$ws.Range("M22").Value = -952.4166
$ws.Range("N22").Value = -1571.5

# Row 27
$ws.Range("H27").Value = 1180.9375
$ws.Range("I27").Value = 1247.4166
$ws.Range("J27").Value = 981.5
$ws.Range("K27").Value = 1247.4166
$ws.Range("L27").Value = 981.5
$ws.Range("M27").Value = -1140.4166
$ws.Range("N27").Value = -1195.5

# Row 46
$ws.Range("H46").Value = 2956.5625
$ws.Range("I46").Value = 2139.7273
$ws.Range("J46").Value = 4753.6
$ws.Range("K46").Value = 2139.7273
$ws.Range("L46").Value = 4753.6
$ws.Range("M46").Value = -1951.7273
$ws.Range("N46").Value = -5129.6

# Row 122
$ws.Range("H122").Value = 7606.4375
$ws.Range("J122").Value = 7120.5
$ws.Range("L122").Value = 21361.5
$ws.Range("N122").Value = -26261.5

# Row 136
$ws.Range("H136").Value = 3309.9412
$ws.Range("I136").Value = 2876
$ws.Range("K136").Value = 8628
$ws.Range("M136").Value = -6078

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4887.857
$ws.Range("I132").Value = 4866.5454
$ws.Range("J132").Value = 4966
$ws.Range("K132").Value = 14599.6362
$ws.Range("L132").Value = 14898
$ws.Range("M132").Value = -12069.6362
$ws.Range("N132").Value = -19958
